$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: the old "HW 2" (E25) and "Lab Notebook Part 1" (G25) entries move away ---
$ws.Range("E25").ClearContents()
$ws.Range("G25").ClearContents()

# --- Row 26: NO CLASS / NA ---
$ws.Range("B26").Value = 25
$ws.Range("B26").Font.Bold = $true
$ws.Range("B26").HorizontalAlignment = -4131

$ws.Range("C26").Value = "NO CLASS"
$ws.Range("C26").HorizontalAlignment = -4131

$ws.Range("D26").Value = "NA"

# --- Row 27: 12.1 - Clocking ---
$ws.Range("B27").Value = 26
$ws.Range("B27").Font.Bold = $true
$ws.Range("B27").HorizontalAlignment = -4131

$ws.Range("C27").Value = "12.1 - Clocking"
$ws.Range("C27").HorizontalAlignment = -4131

$ws.Range("D27").Value = 12

$ws.Range("H27").Value = "Missing"

# --- Row 28: 12.2 - Clocking + Final Project Work Day 3 ---
$ws.Range("B28").Value = 27
$ws.Range("B28").Font.Bold = $true
$ws.Range("B28").HorizontalAlignment = -4131

$ws.Range("C28").Value = "12.2 - Clocking + Final Project Work Day 3"
$ws.Range("C28").HorizontalAlignment = -4131

$ws.Range("D28").Value = 12

$ws.Range("H28").Value = "https://iu.zoom.us/rec/share/r40LjSsDuljX-ADcPdAabvBnSKCeUdVrQodbCKYzzWYz52RUyh2j7RLCF__uBrEB.2FuwHK_E1shbLXcE"
$ws.Hyperlinks.Add($ws.Range("H28"), "https://iu.zoom.us/rec/share/r40LjSsDuljX-ADcPdAabvBnSKCeUdVrQodbCKYzzWYz52RUyh2j7RLCF__uBrEB.2FuwHK_E1shbLXcE")
$ws.Range("H28").Style = "Hyperlink"

# --- Row 29: Final Project Work Day 4 ---
$ws.Range("B29").Value = 28
$ws.Range("B29").Font.Bold = $true
$ws.Range("B29").HorizontalAlignment = -4131

$ws.Range("C29").Value = "Final Project Work Day 4"
$ws.Range("C29").HorizontalAlignment = -4131

# --- Row 30: Lab Notebook Part 1 moves here ---
$ws.Range("B30").Value = 29
$ws.Range("B30").Font.Bold = $true
$ws.Range("B30").HorizontalAlignment = -4131

$ws.Range("G30").Value = "Lab Notebook Part 1"

# --- Update selection to match the final state ---
$ws.Range("C34").Select()
